# DecisionRationale document - XML as storage medium
#
# 1. Remove the "_GoBack" bookmark that currently sits between the bold
#    lead-in of the first bullet and the rest of its body text.
# 2. Append a brand-new bulleted paragraph (same list, numId 11) describing
#    the decision to use XML files as a persistent storage medium, and move
#    the "_GoBack" bookmark to the very end of that new paragraph.

$d = $word.ActiveDocument

# --- 1. Remove the stray _GoBack bookmark from the first bullet ---------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. Create a new paragraph at the end of the document, matching the
#        bullet-list formatting of the existing list item --------------
$lastParagraph = $d.Paragraphs.Last
$endOfDoc = $lastParagraph.Range
$endOfDoc.Collapse(0)          # wdCollapseEnd
$endOfDoc.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Style = $d.Styles.Item("List Paragraph")
$newPara.Range.ListFormat.ApplyListTemplateWithLevel(
    $d.Paragraphs.Item(2).Range.ListFormat.ListTemplate, $true, 1, $false, 2, 1)

# The insertion point sits right before the new paragraph's end-of-paragraph
# mark; that mark happens to carry stray Bold formatting inherited from
# upstream edits, so every plain run below must explicitly clear Bold again
# after the text has been typed in.
$paraStart = $newPara.Range.Start

$boldLeadIn  = "Use of XML files as a persistent storage medium. "
$bodyPart1   = "The decision was made to use XML files to store persistent data, as unlike many other data storage options, such as JSON, XML does not require 3"
$superscript = "rd"
$bodyPart2   = " party libraries, as an XML library is included with Java. As it is plain text, XML data is easily editable for configuration purposes. "

$fullText = $boldLeadIn + $bodyPart1 + $superscript + $bodyPart2

$insertRange = $d.Range($paraStart, $paraStart)
$insertRange.InsertAfter($fullText)

# Compute the character offsets of each run within the newly inserted text.
$boldStart  = $paraStart
$boldEnd    = $boldStart + $boldLeadIn.Length
$body1Start = $boldEnd
$body1End   = $body1Start + $bodyPart1.Length
$supStart   = $body1End
$supEnd     = $supStart + $superscript.Length
$body2Start = $supEnd
$body2End   = $body2Start + $bodyPart2.Length

# Run 1: bold lead-in.
$boldRange = $d.Range($boldStart, $boldEnd)
$boldRange.Font.Bold = $true

# Run 2: plain body text (explicitly clear the inherited Bold taint).
$body1Range = $d.Range($body1Start, $body1End)
$body1Range.Font.Bold = $false

# Run 3: superscript "rd".
$supRange = $d.Range($supStart, $supEnd)
$supRange.Font.Bold = $false
$supRange.Font.Superscript = $true

# Run 4: remainder of the plain body text.
$body2Range = $d.Range($body2Start, $body2End)
$body2Range.Font.Bold = $false
$body2Range.Font.Superscript = $false

# --- 3. Re-add the _GoBack bookmark at the very end of the new paragraph
$bookmarkRange = $d.Range($body2End, $body2End)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)
